# Fill in the previously-empty rows 109-117 with the "16-09-2020" day's
# last reading and the new "17-09-2020" day's log entries.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# Reference cells whose existing styles we want to replicate onto the
# newly filled-in cells:
#   B108 -> time-of-day style (s=4, numFmt "h:mm;@")
#   D108 -> blood-glucose reading style (s=6)
#   A88  -> date style (s=11, numFmt "yyyy/m/d")
#   C103 -> feeding-note style (s=1)
$timeFmt = $ws.Range("B108")
$numFmt  = $ws.Range("D108")
$dateFmt = $ws.Range("A88")
$noteFmt = $ws.Range("C103")

function Set-TimeCell($addr, $val) {
    $timeFmt.Copy() | Out-Null
    $ws.Range($addr).PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Range($addr).Value = $val
}

function Set-NumCell($addr, $val) {
    $numFmt.Copy() | Out-Null
    $ws.Range($addr).PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Range($addr).Value = $val
}

function Set-DateCell($addr, $val) {
    $dateFmt.Copy() | Out-Null
    $ws.Range($addr).PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Range($addr).Value = $val
}

function Set-NoteCell($addr, $val) {
    $noteFmt.Copy() | Out-Null
    $ws.Range($addr).PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Range($addr).Value = $val
}

# Row 109 - last glucose reading of the 16-09-2020 day
Set-TimeCell "B109" 0.9958333333333333
Set-NumCell  "D109" 20.7

# Row 110 - start of a new day: 17-09-2020
Set-DateCell "A110" "17-09-2020"
Set-TimeCell "B110" 0.0
Set-NoteCell "C110" "罐头25g"

# Row 111
Set-TimeCell "B111" 0.006944444444444444
Set-NumCell  "E111" 0.6

# Row 112
Set-TimeCell "B112" 0.08333333333333333
Set-NumCell  "D112" 21.6

# Row 113
Set-TimeCell "B113" 0.1736111111111111
Set-NumCell  "D113" 16.9

# Row 114
Set-TimeCell "B114" 0.2708333333333333
Set-NumCell  "D114" 13.0

# Row 115
Set-TimeCell "B115" 0.3645833333333333
Set-NumCell  "D115" 12.9

# Row 116
Set-TimeCell "B116" 0.375
Set-NoteCell "C116" "鸡肉33g"

# Row 117
Set-TimeCell "B117" 0.3819444444444444
Set-NumCell  "E117" 0.6
